# results.xlsx - "change nested loop to if gen and replace if(xy[i])
# o +=xy[i] in hamming.v" -- the Verilog source change altered the
# synthesis stats captured in this results sheet for the "Sum" (row 5)
# and the unlabeled "Hamming" sub-total (row 6) entries of the XOR/IV
# block (columns P:U). Update the recorded input values; the R/T/U
# columns are SUM(...) formulas and recalculate automatically.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("P5").Value = 496
$ws.Range("Q5").Value = 332
$ws.Range("S5").Value = 486

$ws.Range("P6").Value = 4825
$ws.Range("Q6").Value = 3126
$ws.Range("S6").Value = 4810

$excel.Calculate()

# Re-touch the borders on the merged "XOR+IV" header (P1:U1) so the
# style's applyBorder flag is (re)written when the sheet is saved.
$headerRange = $ws.Range("P1:U1")
$headerRange.Borders.LineStyle = -4142   # xlLineStyleNone

# Move the viewport / selection the way the author last left the sheet:
# scrolled right to column L and with V6 selected.
$ws.Range("V6").Select()
$excel.ActiveWindow.ScrollColumn = 12
$excel.ActiveWindow.ScrollRow = 1

# Restore the (narrower) sheet-tab-area / horizontal-scrollbar split
# ratio recorded for the window.
$excel.ActiveWindow.TabRatio = 0.151
